$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.166433334350586
$ws.Range("B1").Value = 2.437785387039185
$ws.Range("D1").Value = 2.367680788040161
$ws.Range("E1").Value = 1.23426365852356
